$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: BICECORP placeholder row ("2") re-identified as entity "1" with
#     refreshed capital-structure figures ---
# Force text storage (not the number 1) for B2, then restore the default
# (unstyled) cell style so no stray formatting is introduced.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1"
$ws.Range("B2").Style = "Normal"

$ws.Range("D2").Value = 0.0288
$ws.Range("E2").Value = 0.007890000000000001
$ws.Range("G2").Value = 0.3754045307443366
$ws.Range("H2").Value = 0.3754045307443366
$ws.Range("I2").Value = 0.343042071197411
$ws.Range("J2").Value = 0.2623970930562652
$ws.Range("K2").Value = 8.73
$ws.Range("L2").Value = 0.2825242718446602
$ws.Range("M2").Value = 7.62
$ws.Range("N2").Value = 0.05090180360721443
$ws.Range("O2").Value = 0.872852233676976
$ws.Range("P2").Value = 7.62
$ws.Range("Q2").Value = 0.05090180360721443
$ws.Range("R2").Value = 0.872852233676976
$ws.Range("U2").Value = 5.47
$ws.Range("V2").Value = 0.03653974615898464
$ws.Range("W2").Value = 0.1746
$ws.Range("X2").Value = 0.01816978031035037
$ws.Range("Y2").Value = 0.1564302196896496
$ws.Range("Z2").Value = 0.6668105308588691
$ws.Range("AA2").Value = 0.1749691449166723
$ws.Range("AB2").Value = 0.01816978031035037
$ws.Range("AC2").Value = 0.1567993646063219
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = -5.47
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = -0.03792553560285655
$ws.Range("AK2").Value = -0.1317120154105466
$ws.Range("AM2").Value = -0.023
$ws.Range("AN2").Value = 0
$ws.Range("AP2").Value = -0.5656670113753878
$ws.Range("AQ2").Value = -460.8695652173913

# --- Row 3: Bolsa de Comercio de Santiago row refreshed with the same
#     recalculated capital-structure figures ---
$ws.Range("D3").Value = 0.0288
$ws.Range("E3").Value = 0.007890000000000001
$ws.Range("G3").Value = 0.3754045307443366
$ws.Range("H3").Value = 0.3754045307443366
$ws.Range("I3").Value = 0.343042071197411
$ws.Range("J3").Value = 0.2623970930562652
$ws.Range("K3").Value = 8.73
$ws.Range("L3").Value = 0.2825242718446602
$ws.Range("M3").Value = 7.62
$ws.Range("N3").Value = 0.05090180360721443
$ws.Range("O3").Value = 0.872852233676976
$ws.Range("P3").Value = 7.62
$ws.Range("Q3").Value = 0.05090180360721443
$ws.Range("R3").Value = 0.872852233676976
$ws.Range("U3").Value = 5.47
$ws.Range("V3").Value = 0.03653974615898464
$ws.Range("W3").Value = 0.1746
$ws.Range("X3").Value = 0.01816978031035037
$ws.Range("Y3").Value = 0.1564302196896496
$ws.Range("Z3").Value = 0.6668105308588691
$ws.Range("AA3").Value = 0.1749691449166723
$ws.Range("AB3").Value = 0.01816978031035037
$ws.Range("AC3").Value = 0.1567993646063219
$ws.Range("AG3").Value = -5.47
$ws.Range("AJ3").Value = -0.03792553560285655
$ws.Range("AK3").Value = -0.1317120154105466
$ws.Range("AM3").Value = -0.023
$ws.Range("AP3").Value = -0.5656670113753878
$ws.Range("AQ3").Value = -460.8695652173913

# --- Row 4 (BICECORP S.A.) removed entirely from the database ---
$ws.Rows(4).Delete()
